# Update "想去人数" (interest count) figures across the four sheets of the
# 北京-漫展信息 workbook, per the refreshed data snapshot (456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value  = 13268
$ws.Range("F9").Value  = 126
$ws.Range("F11").Value = 5183
$ws.Range("F13").Value = 3624
$ws.Range("F25").Value = 4774
$ws.Range("F29").Value = 279
$ws.Range("F30").Value = 7182
$ws.Range("F34").Value = 2074
$ws.Range("F36").Value = 124
$ws.Range("F39").Value = 9
$ws.Range("F43").Value = 13
$ws.Range("F45").Value = 1259
$ws.Range("F46").Value = 1888
$ws.Range("F47").Value = 85

# --- 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 930

# --- 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 672
$ws.Range("F4").Value = 49

# --- 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 672
$ws.Range("F8").Value  = 13269
$ws.Range("F11").Value = 5183
$ws.Range("F12").Value = 3624
$ws.Range("F24").Value = 4774
$ws.Range("F28").Value = 279
$ws.Range("F29").Value = 7182
$ws.Range("F34").Value = 2074
$ws.Range("F36").Value = 124
$ws.Range("F38").Value = 9
$ws.Range("F44").Value = 1259
$ws.Range("F45").Value = 1888
$ws.Range("F46").Value = 85
